$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set B6 to the same text already used in B5 ("bigint[]"), matching shared string reuse
$ws.Range("B6").Value = "bigint[]"

# Scroll the frozen-pane view back to the top (C3) and select C6, instead of C18
$ws.Range("C6").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 3
